# Update the cryptocurrency price/volume table (GitHub Actions scheduled refresh).
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h)
# Note: Price values that look like plain numbers (e.g. "1.007") are written with a
# leading apostrophe so Excel stores/keeps them as literal text (matching the workbook's
# original inline-string/text formatting) instead of re-interpreting them as numbers,
# which would silently drop formatting such as trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '29.349.32'
$ws.Cells.Item(2, 5).Value = '  +0.00%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.875.49'
$ws.Cells.Item(3, 5).Value = '  -0.14%  '

# Row 4
$ws.Cells.Item(4, 4).Value = "'1.007"
$ws.Cells.Item(4, 5).Value = '  +0.66%  '

# Row 5
$ws.Cells.Item(5, 4).Value = "'0.7091"
$ws.Cells.Item(5, 5).Value = '  -0.74%  '

# Row 6
$ws.Cells.Item(6, 4).Value = "'242.80"
$ws.Cells.Item(6, 5).Value = '  +0.22%  '

# Row 7
$ws.Cells.Item(7, 4).Value = "'1.005"
$ws.Cells.Item(7, 5).Value = '  +0.46%  '

# Row 8
$ws.Cells.Item(8, 4).Value = "'0.07837"
$ws.Cells.Item(8, 5).Value = '  -3.26%  '

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.3109"
$ws.Cells.Item(9, 5).Value = '  -0.65%  '

# Row 10
$ws.Cells.Item(10, 4).Value = "'24.32"
$ws.Cells.Item(10, 5).Value = '  -3.67%  '

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.08064"
$ws.Cells.Item(11, 5).Value = '  -3.52%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '1.890.10'
$ws.Cells.Item(12, 5).Value = '  +1.04%  '

# Row 13
$ws.Cells.Item(13, 4).Value = "'93.33"
$ws.Cells.Item(13, 5).Value = '  +2.06%  '

# Row 14
$ws.Cells.Item(14, 4).Value = "'5.165"
$ws.Cells.Item(14, 5).Value = '  -1.72%  '

# Row 15
$ws.Cells.Item(15, 4).Value = "'0.6997"
$ws.Cells.Item(15, 5).Value = '  -2.69%  '

# Row 16
$ws.Cells.Item(16, 4).Value = "'6.359"
$ws.Cells.Item(16, 5).Value = '  +1.71%  '

# Row 17
$ws.Cells.Item(17, 2).Value = 'WrappedBTC'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(17, 4).Value = '29.484.57'
$ws.Cells.Item(17, 5).Value = '  +0.47%  '

# Row 18
$ws.Cells.Item(18, 2).Value = 'ShibaInu'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(18, 4).Value = "'0.000008312"
$ws.Cells.Item(18, 5).Value = '  -1.29%  '

# Row 19
$ws.Cells.Item(19, 4).Value = "'250.97"
$ws.Cells.Item(19, 5).Value = '  +4.24%  '

# Row 20
$ws.Cells.Item(20, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(20, 4).Value = '2.141.64'
$ws.Cells.Item(20, 5).Value = '  +0.95%  '

# Row 21
$ws.Cells.Item(21, 2).Value = 'Avalanche'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(21, 4).Value = "'13.20"
$ws.Cells.Item(21, 5).Value = '  -0.39%  '

# Row 22
$ws.Cells.Item(22, 4).Value = "'1.003"
$ws.Cells.Item(22, 5).Value = '  +0.42%  '

# Row 23
$ws.Cells.Item(23, 4).Value = "'7.575"
$ws.Cells.Item(23, 5).Value = '  -2.99%  '

# Row 24
$ws.Cells.Item(24, 4).Value = "'1.009"
$ws.Cells.Item(24, 5).Value = '  +0.79%  '

# Row 25
$ws.Cells.Item(25, 4).Value = "'0.1565"
$ws.Cells.Item(25, 5).Value = '  -1.60%  '

# Row 26
$ws.Cells.Item(26, 4).Value = "'8.983"
$ws.Cells.Item(26, 5).Value = '  -0.96%  '

# Row 27
$ws.Cells.Item(27, 4).Value = "'161.32"
$ws.Cells.Item(27, 5).Value = '  -1.18%  '

# Row 28
$ws.Cells.Item(28, 4).Value = "'18.73"
$ws.Cells.Item(28, 5).Value = '  +0.95%  '

# Row 29
$ws.Cells.Item(29, 4).Value = "'1.506"
$ws.Cells.Item(29, 5).Value = '  +0.03%  '

# Row 30
$ws.Cells.Item(30, 4).Value = "'4.339"
$ws.Cells.Item(30, 5).Value = '  -1.90%  '

# Row 31
$ws.Cells.Item(31, 4).Value = "'4.272"
$ws.Cells.Item(31, 5).Value = '  -1.49%  '

# Row 32
$ws.Cells.Item(32, 4).Value = "'1.230"
$ws.Cells.Item(32, 5).Value = '  +2.22%  '

# Row 33
$ws.Cells.Item(33, 4).Value = "'0.05264"
$ws.Cells.Item(33, 5).Value = '  -2.08%  '

# Row 34
$ws.Cells.Item(34, 4).Value = "'1.904"
$ws.Cells.Item(34, 5).Value = '  -2.44%  '

# Row 35
$ws.Cells.Item(35, 4).Value = "'0.7505"
$ws.Cells.Item(35, 5).Value = '  -0.21%  '

# Row 36
$ws.Cells.Item(36, 4).Value = "'1.166"
$ws.Cells.Item(36, 5).Value = '  -1.12%  '

# Row 37
$ws.Cells.Item(37, 4).Value = "'2.717"
$ws.Cells.Item(37, 5).Value = '  +0.59%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'VeChain'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(38, 4).Value = "'0.01866"
$ws.Cells.Item(38, 5).Value = '  -0.78%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'Maker'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(39, 4).Value = '1.271.79'
$ws.Cells.Item(39, 5).Value = '  -0.58%  '

# Row 40
$ws.Cells.Item(40, 4).Value = "'2.770"
$ws.Cells.Item(40, 5).Value = '  +1.11%  '

# Row 41
$ws.Cells.Item(41, 4).Value = "'6.281"
$ws.Cells.Item(41, 5).Value = '  -4.70%  '

# Row 42
$ws.Cells.Item(42, 4).Value = "'0.9031"
$ws.Cells.Item(42, 5).Value = '  +1.29%  '

# Row 43
$ws.Cells.Item(43, 4).Value = "'111.01"
$ws.Cells.Item(43, 5).Value = '  +0.56%  '

# Row 44
$ws.Cells.Item(44, 4).Value = "'71.58"
$ws.Cells.Item(44, 5).Value = '  -2.12%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'PaxDollar'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(45, 4).Value = "'1.003"
$ws.Cells.Item(45, 5).Value = '  +0.24%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(46, 4).Value = '2.035.23'
$ws.Cells.Item(46, 5).Value = '  +0.69%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(47, 4).Value = "'0.00000000125"
$ws.Cells.Item(47, 5).Value = '  -4.23%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'Mantle'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(48, 4).Value = "'0.5230"
$ws.Cells.Item(48, 5).Value = '  +0.50%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'RenderToken'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(49, 4).Value = "'1.788"
$ws.Cells.Item(49, 5).Value = '  -0.76%  '

# Row 50
$ws.Cells.Item(50, 4).Value = "'9.395"
$ws.Cells.Item(50, 5).Value = '  -0.80%  '

# Row 51
$ws.Cells.Item(51, 4).Value = "'0.4293"
$ws.Cells.Item(51, 5).Value = '  -1.64%  '
